{"js": "const replacements = [\n  [\"2024-03-17 Sunday\", \"2024-03-18 Monday\"],\n  [\"98\u00d769=6762\", \"68\u00d793=6324\"],\n  [\"80\u00d783=6640\", \"22\u00d797=2134\"],\n  [\"40\u00d794=3760\", \"51\u00d761=3111\"],\n  [\"96\u00d795=9120\", \"19\u00d744=836\"],\n  [\"29\u00d797=2813\", \"73\u00d749=3577\"],\n  [\"15\u00d714=210\", \"62\u00d745=2790\"],\n  [\"33\u00d716=528\", \"96\u00d799=9504\"],\n  [\"40\u00d760=2400\", \"60\u00d735=2100\"],\n  [\"28\u00d725=700\", \"31\u00d761=1891\"],\n  [\"11\u00d788=968\", \"34\u00d777=2618\"],\n  [\"80\u00d723=1840\", \"24\u00d757=1368\"],\n  [\"21\u00d727=567\", \"46\u00d735=1610\"],\n  [\"45\u00d748=2160\", \"50\u00d763=3150\"],\n  [\"39\u00d715=585\", \"41\u00d792=3772\"],\n  [\"18\u00d714=252\", \"25\u00d756=1400\"],\n  [\"80\u00d739=3120\", \"73\u00d729=2117\"],\n  [\"73\u00d766=4818\", \"31\u00d790=2790\"],\n  [\"73\u00d763=4599\", \"90\u00d727=2430\"],\n  [\"92\u00d720=1840\", \"25\u00d713=325\"],\n  [\"36\u00d798=3528\", \"18\u00d788=1584\"],\n  [\"82\u00d752=4264\", \"85\u00d791=7735\"],\n  [\"68\u00d723=1564\", \"25\u00d787=2175\"],\n  [\"49\u00d745=2205\", \"63\u00d786=5418\"],\n  [\"59\u00d732=1888\", \"69\u00d757=3933\"],\n  [\"61\u00d752=3172\", \"61\u00d783=5063\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-17 Sunday\", \"2024-03-18 Monday\"),\n    @(\"98\u00d769=6762\", \"68\u00d793=6324\"),\n    @(\"80\u00d783=6640\", \"22\u00d797=2134\"),\n    @(\"40\u00d794=3760\", \"51\u00d761=3111\"),\n    @(\"96\u00d795=9120\", \"19\u00d744=836\"),\n    @(\"29\u00d797=2813\", \"73\u00d749=3577\"),\n    @(\"15\u00d714=210\", \"62\u00d745=2790\"),\n    @(\"33\u00d716=528\", \"96\u00d799=9504\"),\n    @(\"40\u00d760=2400\", \"60\u00d735=2100\"),\n    @(\"28\u00d725=700\", \"31\u00d761=1891\"),\n    @(\"11\u00d788=968\", \"34\u00d777=2618\"),\n    @(\"80\u00d723=1840\", \"24\u00d757=1368\"),\n    @(\"21\u00d727=567\", \"46\u00d735=1610\"),\n    @(\"45\u00d748=2160\", \"50\u00d763=3150\"),\n    @(\"39\u00d715=585\", \"41\u00d792=3772\"),\n    @(\"18\u00d714=252\", \"25\u00d756=1400\"),\n    @(\"80\u00d739=3120\", \"73\u00d729=2117\"),\n    @(\"73\u00d766=4818\", \"31\u00d790=2790\"),\n    @(\"73\u00d763=4599\", \"90\u00d727=2430\"),\n    @(\"92\u00d720=1840\", \"25\u00d713=325\"),\n    @(\"36\u00d798=3528\", \"18\u00d788=1584\"),\n    @(\"82\u00d752=4264\", \"85\u00d791=7735\"),\n    @(\"68\u00d723=1564\", \"25\u00d787=2175\"),\n    @(\"49\u00d745=2205\", \"63\u00d786=5418\"),\n    @(\"59\u00d732=1888\", \"69\u00d757=3933\"),\n    @(\"61\u00d752=3172\", \"61\u00d783=5063\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
